# Aragon hospital COVID data - add 2020-08-14 and 2020-08-15 rows (20 hospitals each)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the existing tail-block formatting (date fmt on col A, bordered cells,
# "observaciones" fill on col H) onto the 40 new rows, so no new style entries
# are created in styles.xml - mirrors what dragging the fill handle down does.
$ws.Range("A2412:H2432").Copy() | Out-Null
$ws.Range("A2433:H2472").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 2020-08-14
$ws.Range("A2433").Value = 44057
$ws.Range("B2433").Value = "Hospital Universitario Miguel Servet"
$ws.Range("C2433").Value = 160
$ws.Range("D2433").Value = 18
$ws.Range("E2433").Value = "Zaragoza"
$ws.Range("F2433").Value = "Zaragoza"
$ws.Range("G2433").Value = 50297
$ws.Range("H2433").Value = "Fuente Aragón Hoy"
$ws.Range("A2434").Value = 44057
$ws.Range("B2434").Value = "Hospital Clínico Universitario"
$ws.Range("C2434").Value = 158
$ws.Range("D2434").Value = 21
$ws.Range("E2434").Value = "Zaragoza"
$ws.Range("F2434").Value = "Zaragoza"
$ws.Range("G2434").Value = 50297
$ws.Range("H2434").Value = "Fuente Aragón Hoy"
$ws.Range("A2435").Value = 44057
$ws.Range("B2435").Value = "Hospital Royo Villanova"
$ws.Range("C2435").Value = 39
$ws.Range("D2435").Value = 3
$ws.Range("E2435").Value = "Zaragoza"
$ws.Range("F2435").Value = "Zaragoza"
$ws.Range("G2435").Value = 50297
$ws.Range("H2435").Value = "Fuente Aragón Hoy"
$ws.Range("A2436").Value = 44057
$ws.Range("B2436").Value = "Hospital Nuestra Señora de Gracia"
$ws.Range("C2436").Value = 14
$ws.Range("D2436").Value = 2
$ws.Range("E2436").Value = "Zaragoza"
$ws.Range("F2436").Value = "Zaragoza"
$ws.Range("G2436").Value = 50297
$ws.Range("H2436").Value = "Fuente Aragón Hoy"
$ws.Range("A2437").Value = 44057
$ws.Range("B2437").Value = "Hospital General de la Defensa"
$ws.Range("C2437").Value = 21
$ws.Range("D2437").Value = 1
$ws.Range("E2437").Value = "Zaragoza"
$ws.Range("F2437").Value = "Zaragoza"
$ws.Range("G2437").Value = 50297
$ws.Range("H2437").Value = "Fuente Aragón Hoy"
$ws.Range("A2438").Value = 44057
$ws.Range("B2438").Value = "Hospital Obispo Polanco"
$ws.Range("C2438").Value = 42
$ws.Range("D2438").Value = 2
$ws.Range("E2438").Value = "Teruel"
$ws.Range("F2438").Value = "Teruel"
$ws.Range("G2438").Value = 44216
$ws.Range("H2438").Value = "Fuente Aragón Hoy"
$ws.Range("A2439").Value = 44057
$ws.Range("B2439").Value = "Hospital de Alcañiz"
$ws.Range("C2439").Value = 9
$ws.Range("E2439").Value = "Alcañiz"
$ws.Range("F2439").Value = "Teruel"
$ws.Range("G2439").Value = 44216
$ws.Range("H2439").Value = "Fuente Aragón Hoy"
$ws.Range("A2440").Value = 44057
$ws.Range("B2440").Value = "Hospital de Barbastro"
$ws.Range("C2440").Value = 30
$ws.Range("D2440").Value = 4
$ws.Range("E2440").Value = "Barbastro"
$ws.Range("F2440").Value = "Huesca"
$ws.Range("G2440").Value = 22125
$ws.Range("H2440").Value = "Fuente Aragón Hoy"
$ws.Range("A2441").Value = 44057
$ws.Range("B2441").Value = "Hospital San Jorge"
$ws.Range("C2441").Value = 31
$ws.Range("D2441").Value = 3
$ws.Range("E2441").Value = "Huesca"
$ws.Range("F2441").Value = "Huesca"
$ws.Range("G2441").Value = 22125
$ws.Range("H2441").Value = "Fuente Aragón Hoy"
$ws.Range("A2442").Value = 44057
$ws.Range("B2442").Value = "Hospital Sagrado Corazón"
$ws.Range("E2442").Value = "Huesca"
$ws.Range("F2442").Value = "Huesca"
$ws.Range("G2442").Value = 22125
$ws.Range("H2442").Value = "Fuente Aragón Hoy"
$ws.Range("A2443").Value = 44057
$ws.Range("B2443").Value = "Hospital Ernest Lluch"
$ws.Range("C2443").Value = 6
$ws.Range("E2443").Value = "Calatayud"
$ws.Range("F2443").Value = "Zaragoza"
$ws.Range("G2443").Value = 50297
$ws.Range("H2443").Value = "Fuente Aragón Hoy"
$ws.Range("A2444").Value = 44057
$ws.Range("B2444").Value = "Hospital San José"
$ws.Range("C2444").Value = 1
$ws.Range("E2444").Value = "Teruel"
$ws.Range("F2444").Value = "Teruel"
$ws.Range("G2444").Value = 44216
$ws.Range("H2444").Value = "Fuente Aragón Hoy"
$ws.Range("A2445").Value = 44057
$ws.Range("B2445").Value = "Hospital Ejea – Cinco Villas"
$ws.Range("C2445").Value = 1
$ws.Range("E2445").Value = "Ejea de los Caballeros"
$ws.Range("F2445").Value = "Zaragoza"
$ws.Range("G2445").Value = 50297
$ws.Range("H2445").Value = "Fuente Aragón Hoy"
$ws.Range("A2446").Value = 44057
$ws.Range("B2446").Value = "MAZ"
$ws.Range("C2446").Value = 4
$ws.Range("D2446").Value = 1
$ws.Range("E2446").Value = "Zaragoza"
$ws.Range("F2446").Value = "Zaragoza"
$ws.Range("G2446").Value = 50297
$ws.Range("H2446").Value = "Fuente Aragón Hoy"
$ws.Range("A2447").Value = 44057
$ws.Range("B2447").Value = "Hospital Viamed Montecanal"
$ws.Range("E2447").Value = "Zaragoza"
$ws.Range("F2447").Value = "Zaragoza"
$ws.Range("G2447").Value = 50297
$ws.Range("H2447").Value = "Fuente Aragón Hoy"
$ws.Range("A2448").Value = 44057
$ws.Range("B2448").Value = "Clínica Montpellier"
$ws.Range("C2448").Value = 7
$ws.Range("E2448").Value = "Zaragoza"
$ws.Range("F2448").Value = "Zaragoza"
$ws.Range("G2448").Value = 50297
$ws.Range("H2448").Value = "Fuente Aragón Hoy"
$ws.Range("A2449").Value = 44057
$ws.Range("B2449").Value = "Hospital Quirón"
$ws.Range("C2449").Value = 10
$ws.Range("D2449").Value = 2
$ws.Range("E2449").Value = "Zaragoza"
$ws.Range("F2449").Value = "Zaragoza"
$ws.Range("G2449").Value = 50297
$ws.Range("H2449").Value = "Fuente Aragón Hoy"
$ws.Range("A2450").Value = 44057
$ws.Range("B2450").Value = "Hospital San Juan de Dios de Zaragoza"
$ws.Range("C2450").Value = 26
$ws.Range("E2450").Value = "Zaragoza"
$ws.Range("F2450").Value = "Zaragoza"
$ws.Range("G2450").Value = 50297
$ws.Range("H2450").Value = "Fuente Aragón Hoy"
$ws.Range("A2451").Value = 44057
$ws.Range("B2451").Value = "Clínica Viamed Santiago"
$ws.Range("E2451").Value = "Huesca"
$ws.Range("F2451").Value = "Huesca"
$ws.Range("G2451").Value = 22125
$ws.Range("H2451").Value = "Fuente Aragón Hoy"
$ws.Range("A2452").Value = 44057
$ws.Range("B2452").Value = "Clínica El Pilar"
$ws.Range("E2452").Value = "Zaragoza"
$ws.Range("F2452").Value = "Zaragoza"
$ws.Range("G2452").Value = 50297
$ws.Range("H2452").Value = "Fuente Aragón Hoy"
# 2020-08-15
$ws.Range("A2453").Value = 44058
$ws.Range("B2453").Value = "Hospital Universitario Miguel Servet"
$ws.Range("C2453").Value = 143
$ws.Range("D2453").Value = 19
$ws.Range("E2453").Value = "Zaragoza"
$ws.Range("F2453").Value = "Zaragoza"
$ws.Range("G2453").Value = 50297
$ws.Range("H2453").Value = "Fuente Aragón Hoy"
$ws.Range("A2454").Value = 44058
$ws.Range("B2454").Value = "Hospital Clínico Universitario"
$ws.Range("C2454").Value = 158
$ws.Range("D2454").Value = 22
$ws.Range("E2454").Value = "Zaragoza"
$ws.Range("F2454").Value = "Zaragoza"
$ws.Range("G2454").Value = 50297
$ws.Range("H2454").Value = "Fuente Aragón Hoy"
$ws.Range("A2455").Value = 44058
$ws.Range("B2455").Value = "Hospital Royo Villanova"
$ws.Range("C2455").Value = 36
$ws.Range("D2455").Value = 3
$ws.Range("E2455").Value = "Zaragoza"
$ws.Range("F2455").Value = "Zaragoza"
$ws.Range("G2455").Value = 50297
$ws.Range("H2455").Value = "Fuente Aragón Hoy"
$ws.Range("A2456").Value = 44058
$ws.Range("B2456").Value = "Hospital Nuestra Señora de Gracia"
$ws.Range("C2456").Value = 13
$ws.Range("D2456").Value = 2
$ws.Range("E2456").Value = "Zaragoza"
$ws.Range("F2456").Value = "Zaragoza"
$ws.Range("G2456").Value = 50297
$ws.Range("H2456").Value = "Fuente Aragón Hoy"
$ws.Range("A2457").Value = 44058
$ws.Range("B2457").Value = "Hospital General de la Defensa"
$ws.Range("C2457").Value = 16
$ws.Range("D2457").Value = 1
$ws.Range("E2457").Value = "Zaragoza"
$ws.Range("F2457").Value = "Zaragoza"
$ws.Range("G2457").Value = 50297
$ws.Range("H2457").Value = "Fuente Aragón Hoy"
$ws.Range("A2458").Value = 44058
$ws.Range("B2458").Value = "Hospital Obispo Polanco"
$ws.Range("C2458").Value = 38
$ws.Range("D2458").Value = 2
$ws.Range("E2458").Value = "Teruel"
$ws.Range("F2458").Value = "Teruel"
$ws.Range("G2458").Value = 44216
$ws.Range("H2458").Value = "Fuente Aragón Hoy"
$ws.Range("A2459").Value = 44058
$ws.Range("B2459").Value = "Hospital de Alcañiz"
$ws.Range("C2459").Value = 4
$ws.Range("E2459").Value = "Alcañiz"
$ws.Range("F2459").Value = "Teruel"
$ws.Range("G2459").Value = 44216
$ws.Range("H2459").Value = "Fuente Aragón Hoy"
$ws.Range("A2460").Value = 44058
$ws.Range("B2460").Value = "Hospital de Barbastro"
$ws.Range("C2460").Value = 23
$ws.Range("D2460").Value = 4
$ws.Range("E2460").Value = "Barbastro"
$ws.Range("F2460").Value = "Huesca"
$ws.Range("G2460").Value = 22125
$ws.Range("H2460").Value = "Fuente Aragón Hoy"
$ws.Range("A2461").Value = 44058
$ws.Range("B2461").Value = "Hospital San Jorge"
$ws.Range("C2461").Value = 27
$ws.Range("D2461").Value = 3
$ws.Range("E2461").Value = "Huesca"
$ws.Range("F2461").Value = "Huesca"
$ws.Range("G2461").Value = 22125
$ws.Range("H2461").Value = "Fuente Aragón Hoy"
$ws.Range("A2462").Value = 44058
$ws.Range("B2462").Value = "Hospital Sagrado Corazón"
$ws.Range("E2462").Value = "Huesca"
$ws.Range("F2462").Value = "Huesca"
$ws.Range("G2462").Value = 22125
$ws.Range("H2462").Value = "Fuente Aragón Hoy"
$ws.Range("A2463").Value = 44058
$ws.Range("B2463").Value = "Hospital Ernest Lluch"
$ws.Range("C2463").Value = 6
$ws.Range("E2463").Value = "Calatayud"
$ws.Range("F2463").Value = "Zaragoza"
$ws.Range("G2463").Value = 50297
$ws.Range("H2463").Value = "Fuente Aragón Hoy"
$ws.Range("A2464").Value = 44058
$ws.Range("B2464").Value = "Hospital San José"
$ws.Range("C2464").Value = 1
$ws.Range("E2464").Value = "Teruel"
$ws.Range("F2464").Value = "Teruel"
$ws.Range("G2464").Value = 44216
$ws.Range("H2464").Value = "Fuente Aragón Hoy"
$ws.Range("A2465").Value = 44058
$ws.Range("B2465").Value = "Hospital Ejea – Cinco Villas"
$ws.Range("C2465").Value = 1
$ws.Range("E2465").Value = "Ejea de los Caballeros"
$ws.Range("F2465").Value = "Zaragoza"
$ws.Range("G2465").Value = 50297
$ws.Range("H2465").Value = "Fuente Aragón Hoy"
$ws.Range("A2466").Value = 44058
$ws.Range("B2466").Value = "MAZ"
$ws.Range("C2466").Value = 4
$ws.Range("D2466").Value = 1
$ws.Range("E2466").Value = "Zaragoza"
$ws.Range("F2466").Value = "Zaragoza"
$ws.Range("G2466").Value = 50297
$ws.Range("H2466").Value = "Fuente Aragón Hoy"
$ws.Range("A2467").Value = 44058
$ws.Range("B2467").Value = "Hospital Viamed Montecanal"
$ws.Range("E2467").Value = "Zaragoza"
$ws.Range("F2467").Value = "Zaragoza"
$ws.Range("G2467").Value = 50297
$ws.Range("H2467").Value = "Fuente Aragón Hoy"
$ws.Range("A2468").Value = 44058
$ws.Range("B2468").Value = "Clínica Montpellier"
$ws.Range("C2468").Value = 6
$ws.Range("E2468").Value = "Zaragoza"
$ws.Range("F2468").Value = "Zaragoza"
$ws.Range("G2468").Value = 50297
$ws.Range("H2468").Value = "Fuente Aragón Hoy"
$ws.Range("A2469").Value = 44058
$ws.Range("B2469").Value = "Hospital Quirón"
$ws.Range("C2469").Value = 10
$ws.Range("D2469").Value = 2
$ws.Range("E2469").Value = "Zaragoza"
$ws.Range("F2469").Value = "Zaragoza"
$ws.Range("G2469").Value = 50297
$ws.Range("H2469").Value = "Fuente Aragón Hoy"
$ws.Range("A2470").Value = 44058
$ws.Range("B2470").Value = "Hospital San Juan de Dios de Zaragoza"
$ws.Range("C2470").Value = 28
$ws.Range("E2470").Value = "Zaragoza"
$ws.Range("F2470").Value = "Zaragoza"
$ws.Range("G2470").Value = 50297
$ws.Range("H2470").Value = "Fuente Aragón Hoy"
$ws.Range("A2471").Value = 44058
$ws.Range("B2471").Value = "Clínica Viamed Santiago"
$ws.Range("E2471").Value = "Huesca"
$ws.Range("F2471").Value = "Huesca"
$ws.Range("G2471").Value = 22125
$ws.Range("H2471").Value = "Fuente Aragón Hoy"
$ws.Range("A2472").Value = 44058
$ws.Range("B2472").Value = "Clínica El Pilar"
$ws.Range("E2472").Value = "Zaragoza"
$ws.Range("F2472").Value = "Zaragoza"
$ws.Range("G2472").Value = 50297
$ws.Range("H2472").Value = "Fuente Aragón Hoy"

# Update the visible extent / selection to match (scroll so row 2445 is the
# top-left cell, then select the last date block like the source workbook).
$ws.Application.Goto($ws.Range("A2445"))
$ws.Range("A2454:A2472").Select() | Out-Null
